# Append the new daily data row (row 73) to Sheet1, mirroring the existing
# rows (date text, weekday text, hour number, ranking number).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 73

# Column A holds a date-looking string ("2025/10/07") that must stay plain
# text (as the other date cells in the sheet are), not get auto-converted
# into a date serial number. Force the cell to Text format before writing
# the value, then clear the formatting again so the cell ends up with no
# explicit style (matching the rest of the sheet's unstyled data rows).
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2025/10/07"
$ws.Cells.Item($newRow, 1).ClearFormats()

$ws.Cells.Item($newRow, 2).Value = "火"
$ws.Cells.Item($newRow, 3).Value = 9
$ws.Cells.Item($newRow, 4).Value = 73
